# Fix failing atomic tabOTTR tests: restore the "data"/"auto" row pairing in
# Sheet1 (A7:C9), nudge the sheet's view (selection, tab ratio, column width)
# back to what the fixture expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- bookViews: tab-split ratio 98% -> 99% ---
$excel.ActiveWindow.TabRatio = 0.99

# --- sheetData: rows 7-9 in column A got shuffled; put the numeric "1"
#     marker back on its own row (7), and swap "auto"/"data" onto rows 8/9 ---
$ws.Range("A7").NumberFormat = "General"
$ws.Range("A7").Value = 1

$ws.Range("A8").Value = "auto"

$ws.Range("A9").NumberFormat = "General"
$ws.Range("A9").Value = "data"

# --- sheetView selection follows the data down to A8 ---
$ws.Range("A8").Select() | Out-Null

# --- default column width trimmed slightly ---
$ws.Columns("A").ColumnWidth = 7.35
